# Auto-committed on 2022/02/14 週一
#
# Update the "製作依據之需求規格書與版本" (spec-doc reference) note that is
# repeated down column M (rows 2-19) of the FT test-case sheet: the URS
# version referenced moves from V1.65 (.docx) to V1.64 (.DOCX), and move the
# active selection from G3 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2:M19").Value = "製作依據之需求規格書與版本：PJ201800012_URS_5管理性作業_V1.64.DOCX"

$ws.Range("A2").Select()
